$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.965.45'
$ws.Range("E2").Value = '  -0.10%  '

$ws.Range("D3").Value = '3.427.04'
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '408.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.63%  '

$ws.Range("E7").Value = '  +5.79%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.737'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.68%  '

$ws.Range("E10").Value = '  +5.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.63'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.54%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.48'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.49%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.141'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.01%  '

$ws.Range("D15").Value = '3.960.77'
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("E16").Value = '  +40.38%  '

$ws.Range("D17").Value = '3.404.30'
$ws.Range("E17").Value = '  -0.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.28%  '

$ws.Range("E19").Value = '  +6.15%  '

$ws.Range("D20").Value = '62.009.58'
$ws.Range("E20").Value = '  +0.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '443.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +41.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.74%  '

$ws.Range("E25").Value = '  +2.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '32.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.16%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.38%  '

$ws.Range("E30").Value = '  -1.78%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.95%  '

$ws.Range("E32").Value = '  -0.62%  '

$ws.Range("E33").Value = '  -0.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '43.03'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.71%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0499'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.33'
$ws.Range("D37").Style = "Normal"

$ws.Range("E38").Value = '  +0.08%  '

$ws.Range("E39").Value = '  +0.22%  '

$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.324'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.39%  '

$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.134'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.94'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.89%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.51%  '

$ws.Range("E45").Value = '  +0.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +13.79%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.62'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.99%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.25'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.143'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +23.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.46%  '

$ws.Range("D51").Value = '3.771.02'
$ws.Range("E51").Value = '  +0.45%  '
